$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.603.82'
$ws.Range("E2").Value = '  -5.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.318.68'
$ws.Range("E3").Value = '  -6.38%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.67'
$ws.Range("E5").Value = '  -4.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '83.22'
$ws.Range("E6").Value = '  -9.91%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.523'
$ws.Range("E7").Value = '  -5.11%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.479'
$ws.Range("E9").Value = '  -6.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0808'
$ws.Range("E10").Value = '  -6.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '29.59'
$ws.Range("E11").Value = '  -10.47%  '
$ws.Range("E12").Value = '  -0.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.678.50'
$ws.Range("E13").Value = '  -6.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.33'
$ws.Range("E14").Value = '  -8.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.63'
$ws.Range("E15").Value = '  -5.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.350.51'
$ws.Range("E16").Value = '  -4.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.737'
$ws.Range("E17").Value = '  -7.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '39.606.81'
$ws.Range("E18").Value = '  -4.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0891'
$ws.Range("E19").Value = '  -5.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.01'
$ws.Range("E20").Value = '  -6.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '67.22'
$ws.Range("E21").Value = '  -4.85%  '
$ws.Range("E22").Value = '  -7.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.31'
$ws.Range("E23").Value = '  -3.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.51'
$ws.Range("E24").Value = '  -8.78%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.77'
$ws.Range("E26").Value = '  -8.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.04'
$ws.Range("E27").Value = '  -8.04%  '
$ws.Range("E28").Value = '  -2.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.09'
$ws.Range("E29").Value = '  -6.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.69'
$ws.Range("E30").Value = '  -8.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '149.70'
$ws.Range("E31").Value = '  -5.00%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.05'
$ws.Range("E33").Value = '  -7.13%  '
$ws.Range("E34").Value = '  -5.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0714'
$ws.Range("E35").Value = '  -6.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.112'
$ws.Range("E36").Value = '  -3.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.72'
$ws.Range("E37").Value = '  -6.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0978'
$ws.Range("E38").Value = '  -5.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.47'
$ws.Range("E39").Value = '  -9.88%  '
$ws.Range("E40").Value = '  -9.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.74'
$ws.Range("E41").Value = '  -7.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.30'
$ws.Range("E42").Value = '  -5.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.931.22'
$ws.Range("E43").Value = '  -3.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0261'
$ws.Range("E44").Value = '  -7.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.32'
$ws.Range("E45").Value = '  -7.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.18'
$ws.Range("E46").Value = '  -3.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.61'
$ws.Range("E47").Value = '  -12.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.560.30'
$ws.Range("E48").Value = '  -6.95%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '91.10'
$ws.Range("E49").Value = '  -6.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.08'
$ws.Range("E50").Value = '  -7.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '62.58'
$ws.Range("E51").Value = '  -7.29%  '
